$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 659.3514
$ws.Range("J17").Value = 550.6462
$ws.Range("L17").Value = 1651.9386
$ws.Range("N17").Value = -1987.9386
$ws.Range("H32").Value = 400
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 400
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 400
$ws.Range("N32").Value = -1052
$ws.Range("M32").ClearContents()
$ws.Range("H40").Value = 2536.5454
$ws.Range("I40").Value = 4000.3333
$ws.Range("J40").Value = 1987.625
$ws.Range("K40").Value = 4000.3333
$ws.Range("L40").Value = 1987.625
$ws.Range("M40").Value = -3825.3333
$ws.Range("N40").Value = -2337.625
$ws.Range("H43").Value = 2137.9
$ws.Range("J43").Value = 2374.75
$ws.Range("L43").Value = 2374.75
$ws.Range("N43").Value = -2512.75
$ws.Range("H88").Value = 11484
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 11780.8
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 11780.8
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -12592.8
$ws.Range("H91").Value = 11484
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 11780.8
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 11780.8
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -14588.8
$ws.Range("H115").Value = 1820.5333
$ws.Range("J115").Value = 2244.889
$ws.Range("L115").Value = 6734.667
$ws.Range("N115").Value = -9868.667000000001
$ws.Range("H116").Value = 725086.1
$ws.Range("I116").Value = 3337666.8
$ws.Range("J116").Value = 12564.182
$ws.Range("K116").Value = 3337666.8
$ws.Range("L116").Value = 12564.182
$ws.Range("M116").Value = -3334224.8
$ws.Range("N116").Value = -19448.182
$ws.Range("H118").Value = 773.4737
$ws.Range("I118").Value = 583.1111
$ws.Range("J118").Value = 944.8
$ws.Range("K118").Value = 1749.3333
$ws.Range("L118").Value = 2834.4
$ws.Range("M118").Value = -92.33329999999978
$ws.Range("N118").Value = -6148.4
$ws.Range("H129").Value = 925.7527
$ws.Range("I129").Value = 480
$ws.Range("J129").Value = 951.0795
$ws.Range("K129").Value = 1440
$ws.Range("L129").Value = 2853.2385
$ws.Range("M129").Value = 3560
$ws.Range("N129").Value = -12853.2385
$ws.Range("H137").Value = 3844.0454
$ws.Range("I137").Value = 2206.9092
$ws.Range("J137").Value = 5481.1816
$ws.Range("K137").Value = 6620.7276
$ws.Range("L137").Value = 16443.5448
$ws.Range("M137").Value = -4070.7276
$ws.Range("N137").Value = -21543.5448
$ws.Range("H138").Value = 5362.03
$ws.Range("I138").Value = 872.5
$ws.Range("J138").Value = 7286.1143
$ws.Range("K138").Value = 2617.5
$ws.Range("L138").Value = 21858.3429
$ws.Range("M138").Value = 2522.5
$ws.Range("N138").Value = -32138.3429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1426.238
$ws.Range("I61").Value = 949.8
$ws.Range("J61").Value = 2617.3333
$ws.Range("K61").Value = 949.8
$ws.Range("L61").Value = 2617.3333
$ws.Range("M61").Value = -737.8
$ws.Range("N61").Value = -3041.3333
$ws.Range("H122").Value = 1760.4333
$ws.Range("I122").Value = 1079.7059
$ws.Range("J122").Value = 2650.6155
$ws.Range("K122").Value = 3239.1177
$ws.Range("L122").Value = 7951.8465
$ws.Range("M122").Value = -789.1176999999998
$ws.Range("N122").Value = -12851.8465
$ws.Range("H132").Value = 3034.842
$ws.Range("I132").Value = 1390.2
$ws.Range("K132").Value = 4170.6
$ws.Range("M132").Value = -1640.6
$ws.Range("H136").Value = 1426.238
$ws.Range("I136").Value = 949.8
$ws.Range("J136").Value = 2617.3333
$ws.Range("K136").Value = 2849.4
$ws.Range("L136").Value = 7851.999899999999
$ws.Range("M136").Value = -299.3999999999996
$ws.Range("N136").Value = -12951.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 58682.855
$ws.Range("J132").Value = 58682.855
$ws.Range("L132").Value = 58682.855
$ws.Range("N132").Value = -68802.85500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 3503.3333
$ws.Range("J8").Value = 3503.3333
$ws.Range("L8").Value = 3503.3333
$ws.Range("N8").Value = -3783.3333
$ws.Range("H25").Value = 33000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 33000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 33000
$ws.Range("N25").Value = -33348
$ws.Range("M25").ClearContents()
$ws.Range("H132").Value = 4203.387
$ws.Range("I132").Value = 3460.4
$ws.Range("J132").Value = 5554.273
$ws.Range("K132").Value = 10381.2
$ws.Range("L132").Value = 16662.819
$ws.Range("M132").Value = -7851.200000000001
$ws.Range("N132").Value = -21722.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4593.6875
$ws.Range("I3").Value = 3547.6155
$ws.Range("J3").Value = 9126.667
$ws.Range("K3").Value = 10642.8465
$ws.Range("L3").Value = 27380.001
$ws.Range("M3").Value = -10530.8465
$ws.Range("N3").Value = -27604.001
$ws.Range("H113").Value = 582.7647
$ws.Range("I113").Value = 590.4762
$ws.Range("J113").Value = 570.3077
$ws.Range("K113").Value = 1771.4286
$ws.Range("L113").Value = 1710.9231
$ws.Range("M113").Value = 398.5714000000003
$ws.Range("N113").Value = -6050.9231
$ws.Range("H139").Value = 2191.5312
$ws.Range("I139").Value = 1306.1428
$ws.Range("J139").Value = 3881.818
$ws.Range("K139").Value = 3918.4284
$ws.Range("L139").Value = 11645.454
$ws.Range("M139").Value = 1221.5716
$ws.Range("N139").Value = -21925.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 85841
$ws.Range("I22").Value = 201338.2
$ws.Range("K22").Value = 201338.2
$ws.Range("M22").Value = -201043.2
$ws.Range("H27").Value = 85841
$ws.Range("I27").Value = 201338.2
$ws.Range("K27").Value = 201338.2
$ws.Range("M27").Value = -201231.2
$ws.Range("H40").Value = 7671.1055
$ws.Range("I40").Value = 7781.6
$ws.Range("J40").Value = 7631.643
$ws.Range("K40").Value = 7781.6
$ws.Range("L40").Value = 7631.643
$ws.Range("M40").Value = -7645.6
$ws.Range("N40").Value = -7903.643
$ws.Range("H46").Value = 1876.65
$ws.Range("I46").Value = 3567
$ws.Range("J46").Value = 1578.3529
$ws.Range("K46").Value = 3567
$ws.Range("L46").Value = 1578.3529
$ws.Range("M46").Value = -3379
$ws.Range("N46").Value = -1954.3529
$ws.Range("H82").Value = 1142.3658
$ws.Range("I82").Value = 539.8571
$ws.Range("J82").Value = 1775
$ws.Range("K82").Value = 539.8571
$ws.Range("L82").Value = 1775
$ws.Range("M82").Value = -178.8570999999999
$ws.Range("N82").Value = -2497
$ws.Range("H85").Value = 1142.3658
$ws.Range("I85").Value = 539.8571
$ws.Range("J85").Value = 1775
$ws.Range("K85").Value = 539.8571
$ws.Range("L85").Value = 1775
$ws.Range("M85").Value = 708.1429
$ws.Range("N85").Value = -4271
$ws.Range("H132").Value = 5346.3076
$ws.Range("I132").Value = 2589
$ws.Range("J132").Value = 6806.0586
$ws.Range("K132").Value = 7767
$ws.Range("L132").Value = 20418.1758
$ws.Range("M132").Value = -5237
$ws.Range("N132").Value = -25478.1758

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11113432
$ws.Range("I132").Value = 1400.2632
$ws.Range("J132").Value = 30306940
$ws.Range("K132").Value = 4200.7896
$ws.Range("L132").Value = 90920820
$ws.Range("M132").Value = -1670.7896
$ws.Range("N132").Value = -90925880
